$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets
#    emergency_type_severity_resourc -> emergency_type_priority_resourc
#    resources                       -> resource
# ---------------------------------------------------------------------------
$wsPriorityResource = $wb.Worksheets.Item("emergency_type_severity_resourc")
$wsPriorityResource.Name = "emergency_type_priority_resourc"

$wsResource = $wb.Worksheets.Item("resources")
$wsResource.Name = "resource"

# ---------------------------------------------------------------------------
# 2. emergency_type_priority_resourc sheet (was emergency_type_severity_resourc)
#    New headers + seeded rows (emergency_type_id x priority_id x resource_type_id
#    -> recommended_quantity)
# ---------------------------------------------------------------------------
$wsPriorityResource.Range("A1").Value = "emergency_type_priority_resource_id"
$wsPriorityResource.Range("B1").Value = "emergency_type_id"
$wsPriorityResource.Range("C1").Value = "priority_id"
$wsPriorityResource.Range("D1").Value = "resource_type_id"
$wsPriorityResource.Range("E1").Value = "recommended_quantity"

$wsPriorityResource.Range("C2").Value = 1
$wsPriorityResource.Range("D2").Value = 1
$wsPriorityResource.Range("E2").Value = 2

$wsPriorityResource.Range("C3").Value = 1
$wsPriorityResource.Range("D3").Value = 2
$wsPriorityResource.Range("E3").Value = 1

$wsPriorityResource.Range("C4").Value = 1
$wsPriorityResource.Range("D4").Value = 3
$wsPriorityResource.Range("E4").Value = 1

$wsPriorityResource.Range("C5").Value = 1
$wsPriorityResource.Range("D5").Value = 4
$wsPriorityResource.Range("E5").Value = 3

# column widths (input is in "characters"; the stored OOXML width carries a
# fixed +5/6 padding on top of whatever we set here)
$wsPriorityResource.Columns.Item(2).ColumnWidth = 16.46484375 - 0.8333333333333334
$wsPriorityResource.Columns.Item(3).ColumnWidth = 9.1328125 - 0.8333333333333334
$wsPriorityResource.Columns.Item(4).ColumnWidth = 14.73046875 - 0.8333333333333334
$wsPriorityResource.Columns.Item(5).ColumnWidth = 20.59765625 - 0.8333333333333334

$wsPriorityResource.Range("D13").Select()

# ---------------------------------------------------------------------------
# 3. resource sheet (was resources) - extend recommended_quantity column with
#    averaged values and apply an integer number format
# ---------------------------------------------------------------------------
$resourceRows = @(
    @(1, 1),
    @(2, 1),
    @(3, 1),
    @(4, 1),
    @(5, 2),
    @(6, 2),
    @(7, 2),
    @(8, 2),
    @(9, 2.3571428571428501),
    @(10, 2.5476190476190399),
    @(11, 2.7380952380952301),
    @(12, 2.9285714285714199),
    @(13, 3.11904761904762),
    @(14, 3.3095238095238102),
    @(15, 3.5),
    @(16, 3.6904761904761898),
    @(17, 3.88095238095238),
    @(18, 4.0714285714285703),
    @(19, 4.2619047619047601),
    @(20, 4.4523809523809499),
    @(21, 1),
    @(22, 1),
    @(23, 1),
    @(24, 2),
    @(25, 2),
    @(26, 2),
    @(27, 2),
    @(28, 2.3571428571428501)
)

foreach ($row in $resourceRows) {
    $r = [int]$row[0] + 1
    $wsResource.Cells.Item($r, 1).Value = $row[0]
    $wsResource.Cells.Item($r, 2).Value = $row[1]
}

$wsResource.Range("B2:B29").NumberFormat = "0"

# ---------------------------------------------------------------------------
# 4. Misc selections to match the saved view state
# ---------------------------------------------------------------------------
$wsLocation = $wb.Worksheets.Item("location")
$wsLocation.Range("A2").Select()

# Activate the "resource" sheet last so it ends up as the active / selected
# tab (workbookView activeTab points at it).
$wsResource.Activate()
$wsResource.Range("A24:A29").Select()
